# ==============================================================
# Update capital structure database for Indonesia Insurance (General)
# Rewrites rows 2-8 of the data grid with refreshed values, inserts a
# new company row, and reorders/replaces existing company rows.
# ==============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous data rows (2-7) so no stale cells remain,
# then rewrite rows 2-8 in full with the refreshed dataset below.
$ws.Range("A2:AQ7").ClearContents()

# Row 2
$ws.Range("A2").Value = "Indonesia"
# B2 holds a numeric-looking label ("6") but must stay a text cell, so force
# it with a leading apostrophe just like the source workbook's inlineStr.
$ws.Range("B2").Value = "'6"
$ws.Range("C2").Value = "Insurance (General)"
$ws.Range("D2").Value = -0.0211
$ws.Range("E2").Value = -0.175
$ws.Range("G2").Value = 0.2055830401777449
$ws.Range("H2").Value = 0.2055830401777449
$ws.Range("I2").Value = 0.05776004443621551
$ws.Range("J2").Value = 0.04568702200596951
$ws.Range("K2").Value = 27.698
$ws.Range("L2").Value = 0.01025661914460285
$ws.Range("M2").Value = 10.369
$ws.Range("N2").Value = 0.001722651311384717
$ws.Range("O2").Value = 0.3743591595061015
$ws.Range("P2").Value = 10.369
$ws.Range("Q2").Value = 0.001722651311384717
$ws.Range("R2").Value = 0.3743591595061015
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 547.7950000000001
$ws.Range("V2").Value = 0.09100779005882834
$ws.Range("W2").Value = 0.001828363604968649
$ws.Range("X2").Value = 0.05268856651270041
$ws.Range("Y2").Value = -0.05086020290773176
$ws.Range("Z2").Value = 1.668572812612569
$ws.Range("AA2").Value = 0.01367479789699813
$ws.Range("AB2").Value = 0.05268414315987453
$ws.Range("AC2").Value = -0.04354723348160756
$ws.Range("AD2").Value = 493.556
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 493.556
$ws.Range("AG2").Value = -54.23900000000009
$ws.Range("AH2").Value = 0.07578285478090262
$ws.Range("AI2").Value = 0.2150855712477819
$ws.Range("AJ2").Value = -0.009092919311761967
$ws.Range("AK2").Value = -0.03104869709273742
$ws.Range("AL2").Value = 123.319
$ws.Range("AM2").Value = 123.319
$ws.Range("AN2").Value = 2.37567507725483
$ws.Range("AO2").Value = 1.264857807799285
$ws.Range("AP2").Value = -0.2610731923332407
$ws.Range("AQ2").Value = 1.264857807799285

# Row 3
$ws.Range("A3").Value = "Indonesia"
$ws.Range("B3").Value = "PT Sinar Mas Multiartha Tbk (IDX:SMMA)"
$ws.Range("C3").Value = "Insurance (General)"
$ws.Range("D3").Value = 0.217
$ws.Range("E3").Value = -0.07829999999999999
$ws.Range("G3").Value = 0.2136835914526563
$ws.Range("H3").Value = 0.2136835914526563
$ws.Range("I3").Value = 0.0617525975298961
$ws.Range("J3").Value = 0.0479222922766567
$ws.Range("K3").Value = 32
$ws.Range("L3").Value = 0.01254655949813762
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 500.1
$ws.Range("V3").Value = 0.09568361841350018
$ws.Range("W3").Value = 0.02188782489740082
$ws.Range("X3").Value = 0.0559420262772192
$ws.Range("Y3").Value = -0.03405420137981838
$ws.Range("Z3").Value = 1.916833260683311
$ws.Range("AA3").Value = 0.09185904376408251
$ws.Range("AB3").Value = 0.0538845749878324
$ws.Range("AC3").Value = 0.03797446877625011
$ws.Range("AD3").Value = 462.2
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 462.2
$ws.Range("AG3").Value = -37.90000000000003
$ws.Range("AH3").Value = 0.08124736324005062
$ws.Range("AI3").Value = 0.2331987891019172
$ws.Range("AJ3").Value = -0.007304334419025966
$ws.Range("AK3").Value = -0.02557527498481681
$ws.Range("AL3").Value = 119.9
$ws.Range("AM3").Value = 119.9
$ws.Range("AN3").Value = 2.233929434509425
$ws.Range("AO3").Value = 1.313594662218515
$ws.Range("AP3").Value = -0.1831802803286613
$ws.Range("AQ3").Value = 1.313594662218515

# Row 4
$ws.Range("A4").Value = "Indonesia"
$ws.Range("B4").Value = "P.T. Asuransi Multi Artha Guna Tbk (IDX:AMAG)"
$ws.Range("C4").Value = "Insurance (General)"
$ws.Range("D4").Value = -0.0211
$ws.Range("E4").Value = -0.175
$ws.Range("G4").Value = 0.1319845857418112
$ws.Range("H4").Value = 0.1319845857418112
$ws.Range("I4").Value = 0.1223506743737958
$ws.Range("J4").Value = 0.1186607334006178
$ws.Range("K4").Value = 6.11
$ws.Range("L4").Value = 0.1177263969171484
$ws.Range("M4").Value = 10.1
$ws.Range("N4").Value = 0.1240786240786241
$ws.Range("O4").Value = 1.653027823240589
$ws.Range("P4").Value = 10.1
$ws.Range("Q4").Value = 0.1240786240786241
$ws.Range("R4").Value = 1.653027823240589
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 30
$ws.Range("V4").Value = 0.3685503685503685
$ws.Range("W4").Value = 0.04482758620689655
$ws.Range("X4").Value = 0.05268117554511614
$ws.Range("Y4").Value = -0.007853589338219592
$ws.Range("Z4").Value = 0.4887005649717513
$ws.Range("AA4").Value = 0.05798956745284427
$ws.Range("AB4").Value = 0.05268117554511614
$ws.Range("AC4").Value = 0.005308391907728131
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -30
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.5836575875486381
$ws.Range("AK4").Value = -0.3080082135523614
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 0
$ws.Range("AP4").Value = -4.219409282700422

# Row 5
$ws.Range("A5").Value = "Indonesia"
$ws.Range("B5").Value = "PT Malacca Trust Wuwungan Insurance Tbk (IDX:MTWI)"
$ws.Range("C5").Value = "Insurance (General)"
$ws.Range("G5").Value = 0.05863402061855671
$ws.Range("H5").Value = 0.05863402061855671
$ws.Range("I5").Value = 0.02461340206185567
$ws.Range("J5").Value = 0.01230670103092784
$ws.Range("K5").Value = 0.005
$ws.Range("L5").Value = 0.0006443298969072165
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 2.31
$ws.Range("V5").Value = 0.3392070484581498
$ws.Range("W5").Value = 0.0005422993492407809
$ws.Range("X5").Value = 0.05268117554511614
$ws.Range("Y5").Value = -0.05213887619587537
$ws.Range("Z5").Value = 1.074792243767313
$ws.Range("AA5").Value = 0.01322714681440443
$ws.Range("AB5").Value = 0.05268117554511614
$ws.Range("AC5").Value = -0.03945402873071172
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = -2.31
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = -0.5133333333333333
$ws.Range("AK5").Value = -0.3432392273402675
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AN5").Value = 0
$ws.Range("AP5").Value = -8.13380281690141

# Row 6
$ws.Range("A6").Value = "Indonesia"
$ws.Range("B6").Value = "PT Equity Development Investment Tbk (IDX:GSMF)"
$ws.Range("C6").Value = "Insurance (General)"
$ws.Range("D6").Value = -0.00426
$ws.Range("E6").Value = -0.4429999999999999
$ws.Range("G6").Value = 0.1366711772665764
$ws.Range("H6").Value = 0.1366711772665764
$ws.Range("I6").Value = 0.04682002706359945
$ws.Range("J6").Value = 0.02341001353179973
$ws.Range("K6").Value = 0.313
$ws.Range("L6").Value = 0.004235453315290933
$ws.Range("M6").Value = 0.269
$ws.Range("N6").Value = 0.005274509803921569
$ws.Range("O6").Value = 0.8594249201277956
$ws.Range("P6").Value = 0.269
$ws.Range("Q6").Value = 0.005274509803921569
$ws.Range("R6").Value = 0.8594249201277956
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 12.5
$ws.Range("V6").Value = 0.2450980392156863
$ws.Range("W6").Value = 0.003114427860696517
$ws.Range("X6").Value = 0.0751670848693347
$ws.Range("Y6").Value = -0.07205265700863818
$ws.Range("Z6").Value = 0.6032653061224491
$ws.Range("AA6").Value = 0.01412244897959184
$ws.Range("AB6").Value = 0.06176288721209524
$ws.Range("AC6").Value = -0.0476404382325034
$ws.Range("AD6").Value = 31.1
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 31.1
$ws.Range("AG6").Value = 18.6
$ws.Range("AH6").Value = 0.3788063337393423
$ws.Range("AI6").Value = 0.2303703703703704
$ws.Range("AJ6").Value = 0.2672413793103449
$ws.Range("AK6").Value = 0.1518367346938776
$ws.Range("AL6").Value = 3.4
$ws.Range("AM6").Value = 3.4
$ws.Range("AN6").Value = 6.835164835164836
$ws.Range("AO6").Value = 1.017647058823529
$ws.Range("AP6").Value = 4.087912087912088
$ws.Range("AQ6").Value = 1.017647058823529

# Row 7
$ws.Range("A7").Value = "Indonesia"
$ws.Range("B7").Value = "PT Asuransi Kresna Mitra Tbk (IDX:ASMI)"
$ws.Range("C7").Value = "Insurance (General)"
$ws.Range("D7").Value = -0.0969
$ws.Range("G7").Value = 0.08922716627634661
$ws.Range("H7").Value = 0.08922716627634661
$ws.Range("I7").Value = -0.823185011709602
$ws.Range("J7").Value = -0.823185011709602
$ws.Range("K7").Value = -6.33
$ws.Range("L7").Value = -0.7412177985948478
$ws.Range("M7").Value = -0
$ws.Range("N7").Value = -0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = -0
$ws.Range("Q7").Value = -0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("U7").Value = 1.98
$ws.Range("V7").Value = 0.003100532414657062
$ws.Range("W7").Value = -0.1652741514360314
$ws.Range("X7").Value = 0.05269595748028467
$ws.Range("Y7").Value = -0.217970108916316
$ws.Range("Z7").Value = 0.2277394063841703
$ws.Range("AA7").Value = -0.187471665911091
$ws.Range("AB7").Value = 0.05268711077463292
$ws.Range("AC7").Value = -0.2401587766857239
$ws.Range("AD7").Value = 0.256
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0.256
$ws.Range("AG7").Value = -1.724
$ws.Range("AH7").Value = 0.0004007162803511276
$ws.Range("AI7").Value = 0.008378059955491556
$ws.Range("AJ7").Value = -0.00270696336492504
$ws.Range("AK7").Value = -0.06033034714445688
$ws.Range("AL7").Value = 0.019
$ws.Range("AM7").Value = 0.019
$ws.Range("AN7").Value = -0.03726346433770015
$ws.Range("AO7").Value = -370
$ws.Range("AP7").Value = 0.2509461426491994
$ws.Range("AQ7").Value = -370

# Row 8
$ws.Range("A8").Value = "Indonesia"
$ws.Range("B8").Value = "PT Asuransi Harta Aman Pratama Tbk (IDX:AHAP)"
$ws.Range("C8").Value = "Insurance (General)"
$ws.Range("D8").Value = -0.144
$ws.Range("G8").Value = -1.011392405063291
$ws.Range("H8").Value = -1.011392405063291
$ws.Range("I8").Value = -0.5683544303797469
$ws.Range("J8").Value = -0.5683544303797469
$ws.Range("K8").Value = -4.4
$ws.Range("L8").Value = -0.5569620253164557
$ws.Range("M8").Value = -0
$ws.Range("N8").Value = -0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = -0
$ws.Range("Q8").Value = -0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("U8").Value = 0.905
$ws.Range("V8").Value = 0.06114864864864865
$ws.Range("W8").Value = -0.2820512820512821
$ws.Range("X8").Value = 0.05268117554511614
$ws.Range("Y8").Value = -0.3347324575963982
$ws.Range("Z8").Value = 0.5467128027681661
$ws.Range("AA8").Value = -0.310726643598616
$ws.Range("AB8").Value = 0.05268117554511614
$ws.Range("AC8").Value = -0.3634078191437321
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = -0.905
$ws.Range("AH8").Value = 0
$ws.Range("AI8").Value = 0
$ws.Range("AJ8").Value = -0.06513134220942784
$ws.Range("AK8").Value = -0.09239407861153651
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0
$ws.Range("AN8").Value = -0
$ws.Range("AP8").Value = 0.2144549763033176
